# Update "想去人数" (column F) values on the "展览" and "全部类型" sheets
# to reflect newly-generated counts, matching gh-pages output at 456a3b4.

$wb = $excel.ActiveWorkbook

# Sheet "展览" - row => new value for column F
$sheet1Updates = @{
    2  = 189
    4  = 151
    5  = 1311
    6  = 18248
    7  = 371
    8  = 263
    9  = 1069
    10 = 6869
    13 = 16
    14 = 113
    15 = 67
    17 = 161
    18 = 1304
    19 = 245
    24 = 34
    26 = 993
    27 = 127
    28 = 5170
    29 = 537
    30 = 40
    32 = 74
    33 = 12104
    37 = 287
    38 = 3926
}

# Sheet "全部类型" - row => new value for column F
$sheet4Updates = @{
    2  = 189
    4  = 151
    5  = 1311
    6  = 18248
    7  = 371
    8  = 263
    9  = 1069
    10 = 6869
    11 = 691
    13 = 16
    14 = 113
    15 = 67
    17 = 161
    18 = 1304
    19 = 245
    24 = 34
    26 = 993
    27 = 127
    28 = 5170
    29 = 537
    32 = 40
    34 = 74
    35 = 12104
    39 = 287
    40 = 3926
}

$ws1 = $wb.Worksheets.Item("展览")
foreach ($row in $sheet1Updates.Keys) {
    $ws1.Cells.Item($row, 6).Value = $sheet1Updates[$row]
}

$ws4 = $wb.Worksheets.Item("全部类型")
foreach ($row in $sheet4Updates.Keys) {
    $ws4.Cells.Item($row, 6).Value = $sheet4Updates[$row]
}
